$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New matchup results entered for spring 2024 week 6 (rows 1433-1453)
$newRows = @(
    @(1433, 3, 12, 2, 8),
    @(1434, 5, 6, 6, 14),
    @(1435, 2, 13, 4, 7),
    @(1436, 5, 6, 2, 14),
    @(1437, 2, 16, 4, 4),
    @(1438, 2, 14, 4, 6),
    @(1439, 4, 14, 5, 6),
    @(1440, 5, 14, 4, 6),
    @(1441, 9, 17, 2, 3),
    @(1442, 2, 18, 3, 2),
    @(1443, 4, 6, 5, 14),
    @(1444, 4, 19, 3, 1),
    @(1445, 7, 5, 8, 15),
    @(1446, 2, 13, 3, 7),
    @(1447, 5, 15, 3, 5),
    @(1448, 3, 16, 4, 4),
    @(1449, 6, 17, 7, 3),
    @(1450, 4, 13, 2, 7),
    @(1451, 5, 19, 6, 1),
    @(1452, 4, 13, 7, 7),
    @(1453, 5, 6, 3, 14)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Update the visible window / active selection to match where the user left off entering data
$ws.Range("H1444").Select() | Out-Null
